# Auto update stock data
# Updates the "Current" column header date (A5) from Oct 22, 2025 to Oct 23, 2025
# across all relevant sheets, and refreshes the updated Altman Z-Score values (B5)
# for Reliance Steel & Aluminum and Kaiser Aluminum.

$wb = $excel.ActiveWorkbook

$oldDate = "Current Oct 25 Oct 22, 2025"
$newDate = "Current Oct 25 Oct 23, 2025"

# Sheets whose "Current" date label needs updating (all except Ryerson Holding)
$dateSheets = @("Alcoa", "Rio Tinto", "Reliance Steel & Aluminum", "Kaiser Aluminum")

foreach ($sheetName in $dateSheets) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("A5").Value = $newDate
}

# Update the Altman Z-Score (B5) values that changed, preserving them as text
# (matching the workbook's existing text-based number storage), using the
# apostrophe text-prefix so the stored value stays a plain string "12.19"/"8.89".
$wsReliance = $wb.Worksheets.Item("Reliance Steel & Aluminum")
$wsReliance.Range("B5").Value = "'12.19"

$wsKaiser = $wb.Worksheets.Item("Kaiser Aluminum")
$wsKaiser.Range("B5").Value = "'8.89"
